$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 11 (pushing existing rows 11-31 down to 13-33)
$ws.Rows("11:12").Insert()

# Row 11: new weekly entry (Especial) for Terminal La Palmera de La Serena - Tuna
$ws.Cells.Item(11, 1).Value = 8
$ws.Cells.Item(11, 2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(11, 3).Value = 'Coquimbo'
$ws.Cells.Item(11, 4).Value = 44687
$ws.Cells.Item(11, 5).Value = 4
$ws.Cells.Item(11, 6).Value = 'Fruta'
$ws.Cells.Item(11, 7).Value = 100107
$ws.Cells.Item(11, 8).Value = 'Otros'
$ws.Cells.Item(11, 9).Value = 100107011
$ws.Cells.Item(11, 10).Value = 'Tuna'
$ws.Cells.Item(11, 11).Value = 'Sin especificar'
$ws.Cells.Item(11, 12).Value = 'Especial'
$ws.Cells.Item(11, 13).Value = 100
$ws.Cells.Item(11, 14).Value = 18000
$ws.Cells.Item(11, 15).Value = 19000
$ws.Cells.Item(11, 16).Value = 18500
$ws.Cells.Item(11, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(11, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(11, 19).Value = 1028
$ws.Cells.Item(11, 20).Value = 18

# Row 12: new weekly entry (Primera) for Terminal La Palmera de La Serena - Tuna
$ws.Cells.Item(12, 1).Value = 8
$ws.Cells.Item(12, 2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(12, 3).Value = 'Coquimbo'
$ws.Cells.Item(12, 4).Value = 44687
$ws.Cells.Item(12, 5).Value = 4
$ws.Cells.Item(12, 6).Value = 'Fruta'
$ws.Cells.Item(12, 7).Value = 100107
$ws.Cells.Item(12, 8).Value = 'Otros'
$ws.Cells.Item(12, 9).Value = 100107011
$ws.Cells.Item(12, 10).Value = 'Tuna'
$ws.Cells.Item(12, 11).Value = 'Sin especificar'
$ws.Cells.Item(12, 12).Value = 'Primera'
$ws.Cells.Item(12, 13).Value = 100
$ws.Cells.Item(12, 14).Value = 14000
$ws.Cells.Item(12, 15).Value = 15000
$ws.Cells.Item(12, 16).Value = 14500
$ws.Cells.Item(12, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(12, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(12, 19).Value = 806
$ws.Cells.Item(12, 20).Value = 18
